$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting existing rows 5-7 down to 6-8.
$ws.Rows.Item(5).Insert()

# New row 5: copy the invariant A-L columns from the (now shifted) row 6,
# then set the weekly-specific values.
$ws.Range("A5").Value2 = $ws.Range("A6").Value2
$ws.Range("B5").Value2 = $ws.Range("B6").Value2
$ws.Range("C5").Value2 = $ws.Range("C6").Value2
$ws.Range("E5").Value2 = $ws.Range("E6").Value2
$ws.Range("F5").Value2 = $ws.Range("F6").Value2
$ws.Range("G5").Value2 = $ws.Range("G6").Value2
$ws.Range("H5").Value2 = $ws.Range("H6").Value2
$ws.Range("I5").Value2 = $ws.Range("I6").Value2
$ws.Range("J5").Value2 = $ws.Range("J6").Value2
$ws.Range("K5").Value2 = $ws.Range("K6").Value2
$ws.Range("L5").Value2 = $ws.Range("L6").Value2

$ws.Range("D5").Value2 = 44540
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat

$ws.Range("M5").Value2 = 240
$ws.Range("N5").Value2 = 3500
$ws.Range("O5").Value2 = 3800
$ws.Range("P5").Value2 = 3650
$ws.Range("Q5").Value2 = '$/bandeja 2 kilos'
$ws.Range("R5").Value2 = 'Región del Maule'
$ws.Range("S5").Value2 = 1825
$ws.Range("T5").Value2 = 2

# Update the date and origin on the row that used to be row 5 (now row 6).
$ws.Range("D6").Value2 = 44539
$ws.Range("R6").Value2 = 'Región del Maule'
